$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Metrics" table (A1:F80) lists one row per metric, grouped by
# "Metric Catergory" (col A) with a numeric sort key (col B), the metric
# name (col C), a numeric sort key (col D) and two "X" marker columns
# (E/F) showing which data source(s) supply that metric.
#
# This change adds three new metric rows to the "Aged Care" section:
#   - "% Aged Care Resident Cases (Weekly)"           (new row 55)
#   - "% Aged Care Staff Cases (Weekly)"               (new row 63)
#   - "% Aged Care Staff Cases (Weekly) per 1M"        (new row 64)
# All existing rows from the insertion points onward shift down
# accordingly, and the col D sort key is simply 10 * (row - 1) for every
# data row, so it gets rewritten in full after the inserts.

# 1) Insert a new row directly below "% Aged Care Resident Cases Weekly
#    Change" (old row 54) for the new weekly resident-cases percentage.
$ws.Rows("55:55").Insert()
$ws.Range("A55").Value2 = "Aged Care"
$ws.Range("B55").Value2 = 60
$ws.Range("C55").Value2 = "% Aged Care Resident Cases (Weekly)"
$ws.Range("F55").Value2 = "X"

# 2) Insert two new rows directly below "% Aged Care Staff Cases Weekly
#    Change" (old row 61, now row 62 after the first insert) for the new
#    weekly staff-cases percentage metrics.
$ws.Rows("63:64").Insert()
$ws.Range("A63").Value2 = "Aged Care"
$ws.Range("B63").Value2 = 60
$ws.Range("C63").Value2 = "% Aged Care Staff Cases (Weekly)"
$ws.Range("F63").Value2 = "X"

$ws.Range("A64").Value2 = "Aged Care"
$ws.Range("B64").Value2 = 60
$ws.Range("C64").Value2 = "% Aged Care Staff Cases (Weekly) per 1M"
$ws.Range("F64").Value2 = "X"

# 3) The "Metric - Sort" column (D) is just 10, 20, 30, ... down the
#    whole table - rewrite it for every data row now that 3 rows were
#    inserted.
$lastRow = 83
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value2 = 10 * ($r - 1)
}

# 4) Grow the "Metrics" table / AutoFilter to cover the 3 new rows.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F" + $lastRow))

# 5) Restore the selection shown in the edited workbook.
$ws.Range("F54:F64").Select()
